# Update the Samit ("סאמיט") import template:
#   - add a new subject "חונכות טיפולית" (therapeutic tutoring) to the
#     פרוייקט (project) list on the טבלאות (Tables) sheet, and
#   - add a new רכז (coordinator) named "שיר" (Sir) to the רכזים list.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ראשי")       # main data-entry sheet
$ws2 = $wb.Worksheets.Item("טבלאות")     # lookup/reference lists sheet

# --- טבלאות!A column: פרוייקט (project/subject) list ---------------------
# A new option, "חונכות טיפולית", is inserted above the existing entries,
# pushing the existing ones down a row (B/C columns are untouched).
$ws2.Range("A3").Value() = "חונכות טיפולית"
$ws2.Range("A4").Value() = "שיעורי עזר"
$ws2.Range("A5").Value() = "הוראה מתקנת"
$ws2.Range("A6").Value() = "תרגום"
$ws2.Range("A7").Value() = "חונכות טיפולית"

# --- טבלאות!C column: רכזים (coordinators) list ---------------------------
# Append the new coordinator name at the next free row.
$ws2.Range("C14").Value() = "שיר"

# --- Selections, matching where the editor ended up in the UI -------------
[void]$ws2.Activate()
[void]$ws2.Range("C3:C14").Select()

[void]$ws1.Activate()
[void]$ws1.Range("L2").Select()
